$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(2, -247.4245882215916),
    @(3, -247.42458795326),
    @(4, -247.4245876849283),
    @(5, -247.4245874165966),
    @(6, -247.4245871482649),
    @(7, -247.4245868799333),
    @(8, -247.4245866116016),
    @(9, -247.4245863432699),
    @(10, -247.4245860749382),
    @(11, -247.4245858066066),
    @(12, -247.4245855382749),
    @(13, -247.4245852699432),
    @(14, -247.4245850016115),
    @(15, -247.4245847332799),
    @(16, -247.4245844649482),
    @(17, -247.4245841966165),
    @(18, -247.4245839282849),
    @(19, -247.4245836599532),
    @(20, -247.4245833916215),
    @(21, -247.4245831232898),
    @(22, -247.4245828549582),
    @(23, -247.4245825866265),
    @(24, -247.4245823182948),
    @(25, -247.4245820499631),
    @(26, -247.4245817816315),
    @(27, -247.4245815132998),
    @(28, -247.4245812449682),
    @(29, -247.4245809766365),
    @(30, -247.4245807083048),
    @(31, -247.4245804399731),
    @(32, -247.4245801716414),
    @(33, -247.4245799033098),
    @(34, -247.4245796349781),
    @(35, -247.4245793666464),
    @(36, -247.4245790983147),
    @(37, -247.4245788299831),
    @(38, -247.4245785616514),
    @(39, -247.4245782933197),
    @(40, -247.4245780249881),
    @(41, -247.4245777566564),
    @(42, -247.4245774883247),
    @(43, -247.4245772199931),
    @(44, -247.4245769516614),
    @(45, -247.4245766833297),
    @(46, -247.424576414998),
    @(47, -247.4245761466663),
    @(48, -247.4245758783347),
    @(49, -247.424575610003),
    @(50, -247.4245753416713),
    @(51, -247.4245750733397),
    @(52, -247.424574805008),
    @(53, -247.4245745366763),
    @(54, -247.4245742683447),
    @(55, -247.424574000013),
    @(56, -247.4245737316813),
    @(57, -247.4245734633496),
    @(58, -247.4245731950179),
    @(59, -247.4245729266863),
    @(60, -247.4245726583546),
    @(61, -247.4245723900229),
    @(62, -247.4245721216913),
    @(63, -247.4245718533596),
    @(64, -247.4245715850279),
    @(65, -247.4245713166962),
    @(66, -247.4245710483646),
    @(67, -247.4245707800329),
    @(68, -247.4245705117012),
    @(69, -247.4245702433695),
    @(70, -247.4245699750379),
    @(71, -247.4245697067062),
    @(72, -247.4245694383745),
    @(73, -247.4245691700428),
    @(74, -247.4245689017112),
    @(75, -247.4245686333795),
    @(76, -247.4245683650478),
    @(77, -247.4245680967161),
    @(78, -247.4245678283845),
    @(79, -247.4245675600528),
    @(80, -247.4245672917211),
    @(81, -247.4245670233894),
    @(82, -247.4245667550578),
    @(83, -247.4245664867261),
    @(84, -247.4245662183944),
    @(85, -247.4245659500628),
    @(86, -247.4245656817311),
    @(87, -247.4245654133994),
    @(88, -247.4245651450677),
    @(89, -247.4245648767361),
    @(90, -247.4245646084044),
    @(91, -247.4245643400727),
    @(92, -247.424564071741),
    @(93, -247.4245638034094),
    @(94, -247.4245635350777),
    @(95, -247.424563266746),
    @(96, -247.4245629984144),
    @(97, -247.4245627300827),
    @(98, -247.424562461751),
    @(99, -247.4245621934193),
    @(100, -247.4245619250877),
    @(101, -247.424561656756)
)

foreach ($item in $values) {
    $row = $item[0]
    $val = $item[1]
    $ws.Cells.Item($row, 2).Value = $val
}
